# Update xlsx file to test
# - Insert two new columns (Ebook, Type) before the Description column
# - Populate the new columns with data (Ebook = 0, Type = "B")
# - Refresh the Date column (B) values to the new 2024-08 dates
# - Move the selection to B3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns at L:M (old column L "Description" shifts to N) ---
$ws.Range("L1:M1").EntireColumn.Insert()

# --- Column widths for the two new columns ---
$ws.Columns(12).ColumnWidth = 9.125
$ws.Columns(13).ColumnWidth = 7.625

# --- Header row (row 1): new "Ebook" / "Type" headers, bold + wrap like the others ---
$ws.Range("L1").Font.Bold = $true
$ws.Range("L1").WrapText = $true
$ws.Range("L1").HorizontalAlignment = -4131
$ws.Range("L1").VerticalAlignment = -4108
$ws.Range("L1").Value = "Ebook"

$ws.Range("M1").Font.Bold = $true
$ws.Range("M1").WrapText = $true
$ws.Range("M1").HorizontalAlignment = -4131
$ws.Range("M1").VerticalAlignment = -4108
$ws.Range("M1").Value = "Type"

# --- L column (Ebook): numeric 0 flag for every data row ---
# NumberFormat is primed with "general" first so the COM bridge writes a true
# number (not text) and resolves to the builtin "General" numFmtId=0.
$ws.Range("L2").NumberFormat = "general"
$ws.Range("L2").HorizontalAlignment = -4131
$ws.Range("L2").VerticalAlignment = -4108
$ws.Range("L2").WrapText = $true
$ws.Range("L2").Value = 0

$ws.Range("L3").NumberFormat = "general"
$ws.Range("L3").HorizontalAlignment = -4131
$ws.Range("L3").VerticalAlignment = -4108
$ws.Range("L3").WrapText = $true
$ws.Range("L3").Value = 0

$ws.Range("L4").NumberFormat = "general"
$ws.Range("L4").HorizontalAlignment = -4131
$ws.Range("L4").VerticalAlignment = -4108
$ws.Range("L4").Value = 0

$ws.Range("L5").NumberFormat = "general"
$ws.Range("L5").HorizontalAlignment = -4131
$ws.Range("L5").VerticalAlignment = -4108
$ws.Range("L5").Value = 0

# --- M column (Type): "B" for every data row ---
$ws.Range("M2").HorizontalAlignment = -4131
$ws.Range("M2").VerticalAlignment = -4108
$ws.Range("M2").WrapText = $true
$ws.Range("M2").Value = "B"

$ws.Range("M3").HorizontalAlignment = -4131
$ws.Range("M3").VerticalAlignment = -4108
$ws.Range("M3").WrapText = $true
$ws.Range("M3").Value = "B"

$ws.Range("M4").HorizontalAlignment = -4131
$ws.Range("M4").VerticalAlignment = -4108
$ws.Range("M4").Value = "B"

$ws.Range("M5").HorizontalAlignment = -4131
$ws.Range("M5").VerticalAlignment = -4108
$ws.Range("M5").Value = "B"

# --- Refresh the Date column (B) values ---
$ws.Range("B2").Value = "2024-08-01"
$ws.Range("B3").Value = "2024-08-25"
$ws.Range("B4").Value = "2024-08-25"
$ws.Range("B5").Value = "2024-08-25"

# --- Move the active selection to B3 ---
$ws.Range("B3").Select()
